# Add data for 2021-11-03: update sheet name, October label, and row 11/12 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name + sheet element name)
$ws.Name = "Through 2021-10-26"

# Update the "October (through 10-25)" label to "October (through 10-26)"
$ws.Range("A11").Value = "October (through 10-26)"

# Update October row (row 11) values for years 2015-2021 (columns B-H)
$ws.Range("B11").Value = 25
$ws.Range("C11").Value = 44
$ws.Range("D11").Value = 59
$ws.Range("E11").Value = 56
$ws.Range("F11").Value = 47
$ws.Range("G11").Value = 126
$ws.Range("H11").Value = 165

# Update Total row (row 12) values for years 2015-2021 (columns B-H)
$ws.Range("B12").Value = 251
$ws.Range("C12").Value = 473
$ws.Range("D12").Value = 686
$ws.Range("E12").Value = 604
$ws.Range("F12").Value = 469
$ws.Range("G12").Value = 1027
$ws.Range("H12").Value = 1412
